$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Tier" values in M2 and M3 (value "Search(T1)" removed entirely)
$ws.Range("M2").Value = ""
$ws.Range("M3").Value = ""

# Touch the conditional formatting (re-apply the "unique values" rule on a scratch
# range and remove it again) which causes Excel to persist a fresh differential
# format entry in styles.xml alongside the pre-existing one used by E2:E3.
$tmpRange = $ws.Range("Z1")
$tmpFcs = $tmpRange.FormatConditions
$tmpFc = $tmpFcs.AddUniqueValues()
$tmpFc.Font.Color = 393372
$tmpFc.Interior.Color = 13551615
$tmpFcs.Delete()
